$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D22").Value = 78.06
$ws.Range("E22").Value = 75.18
$ws.Range("F22").Value = 91.37
$ws.Range("G22").Value = 84.65
$ws.Range("A23").Value = 15

$ws.Range("D22:E22").Select()
